$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D ("Price") updates ---
# These prices are stored as literal text in the sheet (e.g. "536.81", "59.187.61"),
# not real numbers, so a leading apostrophe forces Excel to keep them as text instead
# of re-parsing them as numeric values (which would also eat trailing zeros, e.g. 5.10).
# Resetting the cell style to "Normal" afterwards avoids leaving a text-number-format
# behind, matching the original (unstyled) cells.
$ws.Range("D2").Value = "'59.185.92"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'2.521.21"
$ws.Range("D3").Style = "Normal"
$ws.Range("D5").Value = "'535.32"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'134.94"
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").Value = "'0.566"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Value = "'2.519.66"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Value = "'0.0997"
$ws.Range("D10").Style = "Normal"
$ws.Range("D14").Value = "'2.966.78"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = "'59.138.18"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = "'22.45"
$ws.Range("D16").Style = "Normal"
$ws.Range("D18").Value = "'2.519.52"
$ws.Range("D18").Style = "Normal"
$ws.Range("D21").Value = "'321.72"
$ws.Range("D21").Style = "Normal"
$ws.Range("D24").Value = "'65.81"
$ws.Range("D24").Style = "Normal"
$ws.Range("D26").Value = "'0.997"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Value = "'0.161"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Value = "'7.52"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Value = "'0.0₃0766"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Value = "'172.51"
$ws.Range("D30").Style = "Normal"
$ws.Range("D41").Value = "'0.795"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Value = "'278.64"
$ws.Range("D42").Style = "Normal"
$ws.Range("D44").Value = "'5.10"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = "'132.25"
$ws.Range("D45").Style = "Normal"
$ws.Range("D50").Value = "'17.14"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Value = "'1.760.26"
$ws.Range("D51").Style = "Normal"

# --- Column E ("Volume(1h)") updates ---
$ws.Range("E2").Value = '  +2.54%  '
$ws.Range("E3").Value = '  +3.31%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("E5").Value = '  +5.42%  '
$ws.Range("E6").Value = '  +4.37%  '
$ws.Range("E7").Value = '  +0.20%  '
$ws.Range("E8").Value = '  +2.98%  '
$ws.Range("E10").Value = '  +4.44%  '
$ws.Range("E11").Value = '  -1.47%  '
$ws.Range("E13").Value = '  +0.82%  '
$ws.Range("E14").Value = '  +3.29%  '
$ws.Range("E15").Value = '  +2.59%  '
$ws.Range("E16").Value = '  +2.29%  '
$ws.Range("E17").Value = '  +3.19%  '
$ws.Range("E18").Value = '  +2.95%  '
$ws.Range("E19").Value = '  +1.93%  '
$ws.Range("E20").Value = '  +3.35%  '
$ws.Range("E21").Value = '  +2.05%  '
$ws.Range("E22").Value = '  +8.91%  '
$ws.Range("E23").Value = '  +0.02%  '
$ws.Range("E24").Value = '  +3.81%  '
$ws.Range("E25").Value = '  +1.20%  '
$ws.Range("E26").Value = '  +0.21%  '
$ws.Range("E27").Value = '  +0.60%  '
$ws.Range("E28").Value = '  +3.31%  '
$ws.Range("E29").Value = '  +5.79%  '
$ws.Range("E30").Value = '  +1.54%  '
$ws.Range("E31").Value = '  +5.54%  '
$ws.Range("E32").Value = '  +4.41%  '
$ws.Range("E33").Value = '  +0.58%  '
$ws.Range("E34").Value = '  +0.08%  '
$ws.Range("E35").Value = '  +0.22%  '
$ws.Range("E36").Value = '  +2.63%  '
$ws.Range("E37").Value = '  -0.68%  '
$ws.Range("E38").Value = '  +0.95%  '
$ws.Range("E39").Value = '  +4.38%  '
$ws.Range("E40").Value = '  +1.12%  '
$ws.Range("E41").Value = '  +3.35%  '
$ws.Range("E42").Value = '  +2.03%  '
$ws.Range("E43").Value = '  +2.97%  '
$ws.Range("E44").Value = '  +1.82%  '
$ws.Range("E45").Value = '  +10.03%  '
$ws.Range("E46").Value = '  +2.37%  '
$ws.Range("E47").Value = '  +2.88%  '
$ws.Range("E48").Value = '  +5.57%  '
$ws.Range("E49").Value = '  +4.83%  '
$ws.Range("E50").Value = '  +2.74%  '
$ws.Range("E51").Value = '  +3.06%  '
